$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header D1: "standard_quantity" -> "standard_ ng"
$ws.Range("D1").Value = "standard_ ng"

# Add new column E: header "sample_ul" plus a constant 20 ng/ul sample volume
# for each of the three data rows.
$ws.Range("E1").Value = "sample_ul"
$ws.Range("E2").Value = 20
$ws.Range("E3").Value = 20
$ws.Range("E4").Value = 20

# Match formatting: E1 header bold+centered like D1, E2:E4 centered like D2:D4
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2:E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Size the new column
$ws.Columns.Item(5).ColumnWidth = 9.142857142857142
